$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.683.60"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "2.671.51"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.98"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.01"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  +4.17%  "

$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.92"
$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.399"
$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.56"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000196"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").Value = "3.151.61"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "65.490.68"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").Value = "2.674.23"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  -0.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("E20").Value = "  +2.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.71"
$ws.Range("E21").Value = "  -1.47%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.64"
$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000110"
$ws.Range("E24").Value = "  +5.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.74"
$ws.Range("E25").Value = "  +3.73%  "

$ws.Range("E26").Value = "  -3.64%  "

$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("E28").Value = "  -1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.15"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "542.78"
$ws.Range("E30").Value = "  +2.87%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("E32").Value = "  -1.53%  "

$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.58"
$ws.Range("E34").Value = "  +4.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.47"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "157.92"
$ws.Range("E39").Value = "  -2.34%  "

$ws.Range("E40").Value = "  -1.66%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.65"
$ws.Range("E42").Value = "  +1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.62"
$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.08"
$ws.Range("E44").Value = "  -1.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0615"
$ws.Range("E45").Value = "  +1.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("E46").Value = "  -3.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.25"
$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.646"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0260"
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +2.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.03"
$ws.Range("E51").Value = "  +2.11%  "
